$wb = $excel.ActiveWorkbook

# --- workbook.xml: active tab switches from 'settings' (index 2) to 'survey' (index 1) ---
$wb.Worksheets.Item(2).Activate()

# --- survey sheet (index 2): url formula cells get a "'?' + " prefix prepended ---
$survey = $wb.Worksheets.Item(2)

$survey.Range("B9").Value  = "''?' + opendatakit.getHashString('../tables/breathcounter/forms/breathcounter/',null)"
$survey.Range("B12").Value = "''?' + opendatakit.getHashString('../tables/complex_validate_test/forms/complex_validate_test/',null)"
$survey.Range("B15").Value = "''?' + opendatakit.getHashString('../tables/customAppearance/forms/customAppearance/',null)"
$survey.Range("B18").Value = "''?' + opendatakit.getHashString('../tables/exampleForm/forms/exampleForm/',null)"
$survey.Range("B21").Value = "''?' + opendatakit.getHashString('../tables/household/forms/household/',null)"
$survey.Range("B24").Value = "''?' + opendatakit.getHashString('../tables/household/forms/household_new/',null)"
$survey.Range("B27").Value = "''?' + opendatakit.getHashString('../tables/household_member/forms/household_member/',null)"
$survey.Range("B30").Value = "''?' + opendatakit.getHashString('../tables/imnci/forms/imnci_test/',null)"
$survey.Range("B33").Value = "''?' + opendatakit.getHashString('../tables/refrigerators/forms/refrigerators_init/',null)"
$survey.Range("B36").Value = "''?' + opendatakit.getHashString('../tables/refrigerators/forms/refrigerators_update/',null)"
$survey.Range("B39").Value = "''?' + opendatakit.getHashString('../tables/section_test/forms/section_test/',null)"
$survey.Range("B42").Value = "''?' + opendatakit.getHashString('../tables/selects/forms/selects/',null)"

# A11 label text is unaffected semantically (still "complex_validate_test"); rewritten so the
# shared-string table collapses the old formula-only strings cleanly.
$survey.Range("A11").Value = "complex_validate_test"

# sheetView: scrolled/selected to a different cell, and now the tab is the active one
$survey.Application.ActiveWindow.ScrollRow = 34
$survey.Range("B45").Select()

# --- settings sheet (index 3): shared-string renumbering only (same text) ---
$settings = $wb.Worksheets.Item(3)
$settings.Range("B2").Value = "framework"
$settings.Range("B3").Value = "complex_validate_test"
$settings.Range("C3").Value = "complex_validate_test"

# settings sheet is no longer the active tab
$survey.Activate()
